$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Wins"/"Losses"/"Ties" headers in AD1:AF1, matching the formatting
# (bold, centered, bordered) used by the other header cells in row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row (2-44)
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 75
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
